# Remove warnings from tests
# Update the forecasted (AD/AE, years 2018 & 2019) values on Sheet1 that
# changed as a result of the recalculation used by the tests.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AD4").Value  = 403.1172912525
$ws.Range("AE4").Value  = 430.6893002302324

$ws.Range("AD5").Value  = 401.811680049
$ws.Range("AE5").Value  = 428.1180825066259

$ws.Range("AD7").Value  = 1.354794165210821
$ws.Range("AE7").Value  = 1.458232699724667

$ws.Range("AD12").Value = 1.648416584999993
$ws.Range("AE12").Value = 1.415650699602679

$ws.Range("AD13").Value = -1.305611203500007
$ws.Range("AE13").Value = -2.571217723606537

$ws.Range("AD16").Value = 109.8692829173223
$ws.Range("AE16").Value = 117.1512687879092

$ws.Range("AD17").Value = 24.53359160571212
$ws.Range("AE17").Value = 27.01093567694465

$ws.Range("AD19").Value = 356.4429954781211
$ws.Range("AE19").Value = 377.4615180343461

$ws.Range("AD21").Value = 362.3378954781211
$ws.Range("AE21").Value = 383.3564180343461

$ws.Range("AD23").Value = 764.149575527121
$ws.Range("AE23").Value = 811.474500540972

$ws.Range("AD25").Value = 117.1189770825331
$ws.Range("AE25").Value = 124.5044014876338
